$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-13 Saturday" "2024-07-14 Sunday"

Replace-Text "999×2=" "173×9="
Replace-Text "803×6=" "718×5="
Replace-Text "368×2=" "271×3="
Replace-Text "741×3=" "497×7="
Replace-Text "925×3=" "557×9="
Replace-Text "116×3=" "788×7="
Replace-Text "644×6=" "826×8="
Replace-Text "180×3=" "338×9="
Replace-Text "453×8=" "397×6="
Replace-Text "528×7=" "776×4="
Replace-Text "814×7=" "637×6="
Replace-Text "158×3=" "521×6="
Replace-Text "481×8=" "875×6="
Replace-Text "598×2=" "952×4="
Replace-Text "378×7=" "849×4="
Replace-Text "128×4=" "938×6="
Replace-Text "903×5=" "305×8="
Replace-Text "804×8=" "500×8="
Replace-Text "849×3=" "338×3="
Replace-Text "473×3=" "997×5="
Replace-Text "110×3=" "449×9="
Replace-Text "686×8=" "412×8="
Replace-Text "590×3=" "165×3="
Replace-Text "350×7=" "636×6="
Replace-Text "684×4=" "400×8="
